$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells AD1/AE1/AF1 should look exactly like the existing header
# cells (bold, centered, bordered -> style index 1). Copying an existing
# header cell's formatting onto the new cells reproduces that style without
# introducing any new style table entries, then we overwrite the text.
$ws.Range("AC1").Copy($ws.Range("AD1"))
$ws.Range("AC1").Copy($ws.Range("AE1"))
$ws.Range("AC1").Copy($ws.Range("AF1"))

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record (Wins/Losses/Ties) for every data row (2-52)
for ($row = 2; $row -le 52; $row++) {
    $ws.Cells.Item($row, 30).Value = 97
    $ws.Cells.Item($row, 31).Value = 65
    $ws.Cells.Item($row, 32).Value = 0
}
